# ---------------------------------------------------------------------------
# 544_2001_2009.xlsx rework: the sheet is restructured from
#   NPOZZO | DATA | LIVELLOSTATICOmslm | LIVELLODAPRm
# into an index/codice/data/val layout with two data rows (and the old
# "LIVELLODAPRm" column dropped):
#   (blank) | codice | data       | val
#   0       | 544    | 2004-08-23 | 37.53
#   1       | 544    | 2005-12-22 | 37.5
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: headers shift right by one column; A1 is dropped entirely ---
$ws.Range("A1").Clear()
$ws.Range("B1").Value = "codice"
$ws.Range("C1").Value = "data"
$ws.Range("D1").Value = "val"

# --- Row 2: index 0 / codice 544 / date 2004-08-23 / val 37.53 ---
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 544
$ws.Range("C2").Value = 38222
$ws.Range("C2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D2").Value = 37.53

# --- Row 3: index 1 / codice 544 / date 2005-12-22 / val 37.5 ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 544
$ws.Range("C3").Value = 38708
$ws.Range("C3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D3").Value = 37.5

# --- Re-apply the bold/bordered/centered "index" style to A2:A3 (matches the
#     look already used on the header row, B1:D1) ---
$ws.Range("B1").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Page margins (inches -> points: 0.75in=54pt, 1in=72pt, 0.5in=36pt) ---
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
